# Updated symbol list (crypto prices / volume%) per upstream data refresh.
# Applies targeted Price (col D) and Volume(1h) (col E) cell updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "303.88"
    "E2" = "5.71%"
    "D3" = "32.70"
    "E3" = "12.05%"
    "D4" = "5.298"
    "E4" = "2.31%"
    "D5" = "0.07490"
    "E5" = "7.97%"
    "D6" = "7.824"
    "E6" = "5.79%"
    "D7" = "3.804"
    "E7" = "7.15%"
    "D8" = "1.505"
    "E8" = "7.38%"
    "D9" = "0.9204"
    "E9" = "2.48%"
    "D10" = "0.01751"
    "E10" = "2,612.15%"
    "D11" = "0.1691"
    "E11" = "5.99%"
    "D12" = "0.07802"
    "E12" = "6.89%"
    "D13" = "0.08009"
    "E13" = "4.51%"
    "D14" = "0.03033"
    "E14" = "3.80%"
    "D15" = "0.09903"
    "E15" = "10.21%"
    "D16" = "0.001491"
    "E16" = "-6.87%"
    "D17" = "0.04607"
    "E17" = "1.77%"
    "D18" = "0.006204"
    "E18" = "-4.45%"
    "D19" = "3.474"
    "E19" = "0.59%"
    "D20" = "2.229"
    "E21" = "3.82%"
    "D23" = "4.558"
    "E23" = "13.75%"
    "E24" = "4.26%"
    "D25" = "0.001219"
    "E25" = "0.92%"
    "D26" = "0.004438"
    "E26" = "1.62%"
    "E27" = "19.83%"
    "D28" = "0.0001741"
    "E28" = "7.75%"
    "D40" = "0.04542"
    "E40" = "4.13%"
    "D41" = "0.007173"
    "E41" = "3.32%"
    "D42" = "0.1344"
    "E42" = "8.23%"
    "D43" = "0.002170"
    "E43" = "4.48%"
    "D44" = "0.01264"
    "E44" = "6.43%"
    "D45" = "0.00006155"
    "E45" = "5.89%"
    "D47" = "0.01300"
    "E47" = "-0.41%"
}

foreach ($cellRef in $updates.Keys) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (matching the source data's text-formatted Price/Volume columns)
    # instead of auto-converting to a number or percentage.
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
}
